$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, shifting the existing row 217 (AMM, Amman...)
# and everything below it down by one row.
$ws.Rows("217:217").Insert()

# Populate the newly inserted row 217 with the new colo entry (XAP / Chapeco, Brazil).
$ws.Range("A217").Value = "XAP"
$ws.Range("B217").Value = "Chapeco, Brazil"
$ws.Range("C217").Value = "XAP"
$ws.Range("D217").Value = -27.1341991425
$ws.Range("E217").Value = -52.6566009521
$ws.Range("F217").Value = "BR"
$ws.Range("G217").Value = "South America"
$ws.Range("H217").Value = "Chapeco"

# Match the bold/bordered formatting used for column A identifiers by copying
# the format from the row below (which carries the original style).
$ws.Range("A218").Copy()
$ws.Range("A217").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
